$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Task "crear formulario para cargar parametros del sistema" (row 7) is now finished
$ws.Range("B7").Value = "terminado"

# New tasks for bank movements reporting
$ws.Range("A9").Value = "terminar circuito de movimiento de bancos"
$ws.Range("B9").Value = "en proceso"

$ws.Range("A10").Value = "generar reporte de mov de bancos"
$ws.Range("B10").Value = "en proceso"

$ws.Range("A11").Value = "revisar reporte orden de pago esta fallando"
$ws.Range("B11").Value = "no comenzado"

# Update the saved selection to C13
$ws.Range("C13").Select()
